$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This "Log" sheet has one timestamp appended per run in column A.
# Find the next empty row right after the currently used range and
# append the new log entry there, as plain text (not a date value).
$used = $ws.UsedRange
$nextRow = $used.Row + $used.Rows.Count

$cell = $ws.Cells.Item($nextRow, 1)
$cell.Value = "2025-10-15 12:41:06"
